$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.059.91"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "2.315.73"
$ws.Range("E3").Value = "  -1.65%  "
$ws.Range("D5").Value = "'532.04"
$ws.Range("E5").Value = "  +1.94%  "
$ws.Range("D6").Value = "'132.20"
$ws.Range("E6").Value = "  -3.65%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("D8").Value = "'0.535"
$ws.Range("E8").Value = "  -0.72%  "
$ws.Range("D9").Value = "2.339.76"
$ws.Range("E9").Value = "  -1.31%  "
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("E12").Value = "  -2.94%  "
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("D14").Value = "'23.48"
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("D15").Value = "2.735.08"
$ws.Range("E15").Value = "  -1.58%  "
$ws.Range("D16").Value = "57.120.01"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("E17").Value = "  -2.38%  "
$ws.Range("D18").Value = "2.340.17"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").Value = "'339.45"
$ws.Range("E19").Value = "  +2.93%  "
$ws.Range("E20").Value = "  -2.01%  "
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("D22").Value = "'4.15"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "'61.62"
$ws.Range("E24").Value = "  +0.45%  "
$ws.Range("D25").Value = "'8.76"
$ws.Range("E25").Value = "  +6.22%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").Value = "'0.994"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "'170.47"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").Value = "0.0₃0721"
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("D33").Value = "'18.49"
$ws.Range("E33").Value = "  -0.87%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").Value = "'0.993"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("E37").Value = "  -1.49%  "
$ws.Range("D38").Value = "'0.905"
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "'39.07"
$ws.Range("E40").Value = "  +1.38%  "
$ws.Range("D41").Value = "'148.36"
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("D43").Value = "'3.58"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").Value = "'278.08"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("E45").Value = "  -3.45%  "
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("E49").Value = "  +2.13%  "
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("E51").Value = "  -0.84%  "
